# Rotate the species-record data in rows 2-4 of the active sheet:
#   old row 2 -> new row 4
#   old row 3 -> new row 2
#   old row 4 -> new row 3
# Only the columns whose values actually differ between the three rows
# are touched (A, B, D, E, F, G, H, Q, R, AI, AN, AO); every other
# column is identical across rows 2-4 already, so it is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","Q","R","AI","AN","AO")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value2
    $row3[$col] = $ws.Range($col + "3").Value2
    $row4[$col] = $ws.Range($col + "4").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value2 = $row3[$col]
    $ws.Range($col + "3").Value2 = $row4[$col]
    $ws.Range($col + "4").Value2 = $row2[$col]
}
